$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 22, pushing existing data (rows 22-29) down to rows 24-31.
$ws.Rows("22:23").Insert()

# New row 22 data (week of 2021-08-09 / serial 44417, Calidad "Primera")
$ws.Range("A22").Value = 1
$ws.Range("B22").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C22").Value = "Arica y Parinacota"
$ws.Range("D22").Value = 44417
$ws.Range("D22").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E22").Value = 15
$ws.Range("F22").Value = 100112009
$ws.Range("G22").Value = "Acelga"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 250
$ws.Range("K22").Value = 1800
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = 1900
$ws.Range("N22").Value = "$/atado 2,5 a 3 kilos"
$ws.Range("O22").Value = "Región de Arica y Parinacota"
$ws.Range("P22").Value = 633
$ws.Range("Q22").Value = 3
$ws.Range("R22").Value = "Hortaliza"

# New row 23 data (week of 2021-08-09 / serial 44417, Calidad "Segunda")
$ws.Range("A23").Value = 1
$ws.Range("B23").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C23").Value = "Arica y Parinacota"
$ws.Range("D23").Value = 44417
$ws.Range("D23").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E23").Value = 15
$ws.Range("F23").Value = 100112009
$ws.Range("G23").Value = "Acelga"
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Segunda"
$ws.Range("J23").Value = 200
$ws.Range("K23").Value = 1500
$ws.Range("L23").Value = 1600
$ws.Range("M23").Value = 1550
$ws.Range("N23").Value = "$/atado 2,5 a 3 kilos"
$ws.Range("O23").Value = "Región de Arica y Parinacota"
$ws.Range("P23").Value = 517
$ws.Range("Q23").Value = 3
$ws.Range("R23").Value = "Hortaliza"
